$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for the columns that change, for rows 2..8
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$before = @{}
foreach ($r in 2..8) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowData
}

# New row order: after-row 2 gets before-row 7's data, 3 <- 8, 4 <- 2, 5 <- 3, 6 <- 4, 7 <- 5, 8 <- 6
$mapping = @{ 2 = 7; 3 = 8; 4 = 2; 5 = 3; 6 = 4; 7 = 5; 8 = 6 }

foreach ($destRow in 2..8) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $before[$srcRow][$c]
    }
}
